# Update crypto price/volume data cells per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.843.82'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '3.379.50'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'580.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = "'178.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.31%  '
$ws.Range("D7").Value = "'0.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.11%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '3.376.64'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("D11").Value = "'6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '3.968.38'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").Value = "'28.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.66%  '
$ws.Range("D16").Value = '65.993.32'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '3.379.80'
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("D19").Value = "'5.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").Value = "'13.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").Value = "'364.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = "'7.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").Value = "'72.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("D24").Value = "'0.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").Value = "'0.527"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  +3.90%  '
$ws.Range("D27").Value = "'9.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").Value = "'5.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("D32").Value = "'23.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = "'6.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").Value = "'1.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.33%  '
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").Value = "'162.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.76%  '
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("D39").Value = "'27.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.94%  '
$ws.Range("D40").Value = "'1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").Value = "'2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").Value = '2.678.93'
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("D43").Value = "'4.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("D44").Value = "'6.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("D45").Value = "'0.0679"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").Value = "'24.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("D47").Value = "'39.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").Value = "'330.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.46%  '
$ws.Range("D49").Value = "'0.0283"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").Value = "'31.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.90%  '
